$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.897.81"
$ws.Range("E2").Value = '  +2.65%  '
$ws.Range("D3").Value = "'2.567.13"
$ws.Range("E3").Value = '  +2.23%  '
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = "'600.24"
$ws.Range("E5").Value = '  +1.78%  '
$ws.Range("D6").Value = "'178.33"
$ws.Range("E6").Value = '  +0.53%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").Value = "'0.520"
$ws.Range("E8").Value = '  +0.46%  '
$ws.Range("D9").Value = "'2.565.26"
$ws.Range("E9").Value = '  +2.14%  '
$ws.Range("E10").Value = '  +11.65%  '
$ws.Range("E11").Value = '  +0.04%  '
$ws.Range("E12").Value = '  +1.13%  '
$ws.Range("D13").Value = "'5.02"
$ws.Range("E13").Value = '  +1.52%  '
$ws.Range("D14").Value = "'3.011.43"
$ws.Range("E14").Value = '  +1.09%  '
$ws.Range("E15").Value = '  +5.59%  '
$ws.Range("E16").Value = '  +2.04%  '
$ws.Range("D17").Value = "'69.797.98"
$ws.Range("E17").Value = '  +2.63%  '
$ws.Range("D18").Value = "'2.566.32"
$ws.Range("E18").Value = '  +2.17%  '
$ws.Range("D19").Value = "'7.70"
$ws.Range("E19").Value = '  +2.30%  '
$ws.Range("D20").Value = "'11.20"
$ws.Range("E20").Value = '  +1.66%  '
$ws.Range("D21").Value = "'365.84"
$ws.Range("E21").Value = '  +3.60%  '
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("E23").Value = '  -0.15%  '
$ws.Range("D24").Value = "'70.85"
$ws.Range("E24").Value = '  -0.24%  '
$ws.Range("D25").Value = "'4.29"
$ws.Range("E25").Value = '  -0.61%  '
$ws.Range("D26").Value = "'1.78"
$ws.Range("E26").Value = '  +1.61%  '
$ws.Range("D27").Value = "'9.22"
$ws.Range("E27").Value = '  +0.35%  '
$ws.Range("E28").Value = '  +2.46%  '
$ws.Range("E29").Value = '  +1.09%  '
$ws.Range("D30").Value = "'0.0₃0921"
$ws.Range("E30").Value = '  +0.17%  '
$ws.Range("D31").Value = "'515.67"
$ws.Range("E31").Value = '  +1.28%  '
$ws.Range("D32").Value = "'7.80"
$ws.Range("E32").Value = '  -0.14%  '
$ws.Range("E33").Value = '  +0.85%  '
$ws.Range("E34").Value = '  +1.64%  '
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("D36").Value = "'163.94"
$ws.Range("E36").Value = '  -0.43%  '
$ws.Range("E37").Value = '  -1.75%  '
$ws.Range("D38").Value = "'19.04"
$ws.Range("E38").Value = '  +3.33%  '
$ws.Range("D39").Value = "'18.91"
$ws.Range("E39").Value = '  +1.37%  '
$ws.Range("D40").Value = "'1.36"
$ws.Range("E40").Value = '  +0.98%  '
$ws.Range("D41").Value = "'1.77"
$ws.Range("E41").Value = '  +1.33%  '
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("D43").Value = "'4.93"
$ws.Range("E43").Value = '  +0.78%  '
$ws.Range("E44").Value = '  -1.56%  '
$ws.Range("D45").Value = "'2.47"
$ws.Range("E45").Value = '  -1.08%  '
$ws.Range("D46").Value = "'39.03"
$ws.Range("E46").Value = '  +0.33%  '
$ws.Range("D47").Value = "'152.10"
$ws.Range("E47").Value = '  +3.02%  '
$ws.Range("E48").Value = '  +1.36%  '
$ws.Range("D49").Value = "'0.523"
$ws.Range("E49").Value = '  +0.24%  '
$ws.Range("E50").Value = '  -1.14%  '
$ws.Range("E51").Value = '  +1.55%  '
